# Updates the crypto prices/volumes table, and swaps the Hedera / InternetComputer rows (36 and 37)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "19.997.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -8.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.406.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -8.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "273.71"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3708"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3071"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.05"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -9.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9962"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06573"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -8.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.421"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.162"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -7.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.91"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -9.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.406.52"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001008"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -8.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05748"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -12.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.66"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -11.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.587"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -9.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.42"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.84"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.322"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "19.994.63"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.264"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "138.72"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.90"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.565.47"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -9.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "108.70"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.848"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -20.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.394"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -8.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8465"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -12.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07705"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.418"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05765"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.34%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.819"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1917"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02036"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.23"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.059"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -10.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.268"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -12.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5290"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.525"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5110"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.806"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.61"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.046"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -9.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.17%  "
